$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.556.16"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "1.597.40"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.45"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.35"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.31%  "
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("E10").Value = "  -3.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0871"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").Value = "1.824.54"
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("D13").Value = "1.619.03"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("E14").Value = "  -3.90%  "
$ws.Range("E15").Value = "  -3.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.59"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("D17").Value = "27.568.54"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.03"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.42"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.65%  "
$ws.Range("E20").Value = "  -4.05%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  -2.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.82"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.36"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.70"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("E29").Value = "  -4.57%  "
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("D33").Value = "1.370.31"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.32"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.957"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.33%  "
$ws.Range("E38").Value = "  -2.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.537"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.817"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.58%  "
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.972"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.26%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.32"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.94"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("E45").Value = "  -3.79%  "
$ws.Range("D46").Value = "1.735.01"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.09"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.73"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").Value = "0.0₆01000"
$ws.Range("E49").Value = "  -4.28%  "
$ws.Range("E50").Value = "  -4.30%  "
$ws.Range("E51").Value = "  -1.00%  "
